# Bolnisi Municipality area workbook: drop the old census-years columns
# (1989 / 2002) and the "(according to the population census data)"
# caption row, keeping only the 2014 figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two obsolete year columns (B=1989, C=2002). This shifts the
# old column D (2014 figures) left into column B, which is exactly what
# the target layout needs.
$ws.Range("B:C").Delete()

# Remove the caption row "(according to the population census data)"
# (old row 2), shifting everything below it up by one row.
$ws.Rows("2:2").Delete()

# Match the target row heights (20.1pt, custom height) for the five
# remaining rows.
$ws.Rows("1:5").RowHeight = 20.1
